$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text storage (these columns hold text like "261.61", "1.25%", "3")
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "261.61"
$ws.Range("E2").Value = "1.25%"
$ws.Range("G2").Value = "3"
$ws.Range("D3").Value = "27.49"
$ws.Range("E3").Value = "1.72%"
$ws.Range("G3").Value = "3"
$ws.Range("D4").Value = "4.726"
$ws.Range("E4").Value = "2.83%"
$ws.Range("G4").Value = "3"
$ws.Range("D5").Value = "0.06071"
$ws.Range("G5").Value = "3"
$ws.Range("D6").Value = "6.712"
$ws.Range("E6").Value = "1.04%"
$ws.Range("G6").Value = "3"
$ws.Range("D7").Value = "0.8630"
$ws.Range("E7").Value = "1.12%"
$ws.Range("G7").Value = "3"
$ws.Range("D8").Value = "0.9232"
$ws.Range("E8").Value = "-1.88%"
$ws.Range("G8").Value = "3"
$ws.Range("D9").Value = "0.1407"
$ws.Range("E9").Value = "0.09%"
$ws.Range("G9").Value = "3"
$ws.Range("D10").Value = "0.05134"
$ws.Range("E10").Value = "0.14%"
$ws.Range("G10").Value = "3"
$ws.Range("D11").Value = "0.07098"
$ws.Range("E11").Value = "-0.11%"
$ws.Range("G11").Value = "3"
$ws.Range("D12").Value = "0.03071"
$ws.Range("E12").Value = "-1.38%"
$ws.Range("G12").Value = "3"
$ws.Range("D13").Value = "0.09106"
$ws.Range("E13").Value = "-0.49%"
$ws.Range("G13").Value = "3"
$ws.Range("D14").Value = "0.001542"
$ws.Range("E14").Value = "1.13%"
$ws.Range("G14").Value = "3"
$ws.Range("D15").Value = "0.0006059"
$ws.Range("E15").Value = "-0.18%"
$ws.Range("G15").Value = "3"
$ws.Range("D16").Value = "0.006194"
$ws.Range("E16").Value = "1.22%"
$ws.Range("G16").Value = "3"
$ws.Range("D17").Value = "3.472"
$ws.Range("E17").Value = "-1.18%"
$ws.Range("G17").Value = "3"
$ws.Range("E18").Value = "-0.45%"
$ws.Range("G18").Value = "3"
$ws.Range("E19").Value = "-1.26%"
$ws.Range("G19").Value = "3"
$ws.Range("E20").Value = "2.45%"
$ws.Range("G20").Value = "3"
$ws.Range("D21").Value = "0.1290"
$ws.Range("E21").Value = "0.93%"
$ws.Range("G21").Value = "3"
$ws.Range("D22").Value = "4.097"
$ws.Range("E22").Value = "7.13%"
$ws.Range("G22").Value = "3"
$ws.Range("D23").Value = "0.04268"
$ws.Range("E23").Value = "-0.31%"
$ws.Range("G23").Value = "3"
$ws.Range("D24").Value = "0.001217"
$ws.Range("E24").Value = "-0.25%"
$ws.Range("G24").Value = "3"
$ws.Range("D25").Value = "0.003913"
$ws.Range("E25").Value = "-8.95%"
$ws.Range("G25").Value = "3"
$ws.Range("G26").Value = "3"
$ws.Range("G27").Value = "3"
$ws.Range("G28").Value = "3"
$ws.Range("G29").Value = "3"
$ws.Range("G30").Value = "3"
$ws.Range("G31").Value = "3"
$ws.Range("G32").Value = "3"
$ws.Range("G33").Value = "3"
$ws.Range("G34").Value = "3"
$ws.Range("G35").Value = "3"
$ws.Range("G36").Value = "3"
$ws.Range("G37").Value = "3"
$ws.Range("G38").Value = "3"
$ws.Range("G39").Value = "3"
$ws.Range("D40").Value = "0.03880"
$ws.Range("E40").Value = "1.42%"
$ws.Range("G40").Value = "3"
$ws.Range("E41").Value = "1.25%"
$ws.Range("G41").Value = "3"
$ws.Range("D42").Value = "0.004150"
$ws.Range("E42").Value = "5.50%"
$ws.Range("G42").Value = "3"
$ws.Range("D43").Value = "0.01502"
$ws.Range("E43").Value = "24.71%"
$ws.Range("G43").Value = "3"
$ws.Range("D44").Value = "0.002209"
$ws.Range("E44").Value = "-9.47%"
$ws.Range("G44").Value = "3"
$ws.Range("D45").Value = "0.00005310"
$ws.Range("E45").Value = "-2.78%"
$ws.Range("G45").Value = "3"
$ws.Range("E46").Value = "0.04%"
$ws.Range("G46").Value = "3"
$ws.Range("G47").Value = "3"
$ws.Range("E48").Value = "-47.00%"
$ws.Range("G48").Value = "3"
$ws.Range("E49").Value = "0.04%"
$ws.Range("G49").Value = "3"
$ws.Range("E50").Value = "0.04%"
$ws.Range("G50").Value = "3"
$ws.Range("G51").Value = "3"
